# Installed capacities 2022 Germany
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean: wipe the old PV/OIL/Capacity demo table.
$ws.Cells.ClearContents()

# Header row.
$ws.Range("A1").Value = "Installed Capacities year 2022"
$ws.Range("B1").Value = "MW"

# Data rows (row 2 intentionally left blank, data resumes at row 3).
$data = @(
    @(3,  "Nuclear", 4056),
    @(4,  "Fossil Hard coal", 18830),
    @(5,  "Wind Onshore", 55797),
    @(6,  "Fossil Brown coal/Lignite", 19106),
    @(7,  "Geothermal", 58),
    @(8,  "Hydro Run-of-river and poundage", 3743),
    @(9,  "Hydro Water Reservoir", 1408),
    @(10, "Wind Offshore", 7787),
    @(11, "Hydro Pumped Storage", 9280),
    @(12, "Other renewable", 404),
    @(13, "Solar", 56567),
    @(14, "Waste", 1620),
    @(15, "Fossil Gas", 30649),
    @(16, "Fossil Oil", 3966),
    @(17, "Other", 1679),
    @(18, "Biomass", 8590),
    @(19, "Total Grand capacity", 223540)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
}

# Header cell formatting: bold, 14pt.
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 14
$ws.Rows.Item(1).RowHeight = 19

# Column A width to fit the longest label.
$ws.Columns.Item(1).ColumnWidth = 29.6666666666667

# Selection as left by the author.
$ws.Range("G5").Select() | Out-Null
